$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as TEXT (no numeric auto-conversion),
# using a scratch cell formatted as Text, copied via PasteSpecial values-only,
# then cleared completely so no residue/style is left behind.
$scratch = $ws.Range("ZZ1")
function Set-TextValue([string]$cellAddr, [string]$val) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $scratch.Clear()
}

$ws.Range("D2").Value = "24.898.84"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "1.712.16"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  -0.13%  "
Set-TextValue "D5" "310.91"
$ws.Range("E5").Value = "  +1.33%  "
Set-TextValue "D6" "0.9974"
$ws.Range("E6").Value = "  -0.03%  "
Set-TextValue "D7" "0.3749"
$ws.Range("E7").Value = "  +1.04%  "
Set-TextValue "D8" "49.57"
$ws.Range("E8").Value = "  +2.73%  "
Set-TextValue "D9" "0.3449"
$ws.Range("E9").Value = "  +0.36%  "
Set-TextValue "D10" "1.208"
$ws.Range("E10").Value = "  +2.26%  "
Set-TextValue "D11" "0.07550"
$ws.Range("E11").Value = "  +4.14%  "
Set-TextValue "D12" "0.9992"
$ws.Range("E12").Value = "  -0.14%  "
Set-TextValue "D13" "21.11"
$ws.Range("E13").Value = "  +3.68%  "
Set-TextValue "D14" "6.318"
$ws.Range("E14").Value = "  +3.57%  "
Set-TextValue "D15" "7.047"
$ws.Range("E15").Value = "  +4.49%  "
$ws.Range("D16").Value = "1.711.48"
$ws.Range("E16").Value = "  +1.99%  "
Set-TextValue "D17" "0.00001136"
$ws.Range("E17").Value = "  +2.51%  "
Set-TextValue "D18" "0.06713"
$ws.Range("E18").Value = "  -0.16%  "
Set-TextValue "D19" "0.9975"
$ws.Range("E19").Value = "  -0.09%  "
Set-TextValue "D20" "84.93"
$ws.Range("E20").Value = "  +4.73%  "
$ws.Range("E21").Value = "  +5.53%  "
Set-TextValue "D22" "6.390"
$ws.Range("E22").Value = "  +4.92%  "
Set-TextValue "D23" "13.16"
$ws.Range("E23").Value = "  +10.16%  "
$ws.Range("D24").Value = "24.887.27"
$ws.Range("E24").Value = "  +2.36%  "
Set-TextValue "D25" "2.453"
$ws.Range("E25").Value = "  +1.07%  "
Set-TextValue "D26" "2.799"
$ws.Range("E26").Value = "  +5.32%  "
Set-TextValue "D27" "20.44"
$ws.Range("E27").Value = "  +4.56%  "
Set-TextValue "D28" "151.73"
$ws.Range("E28").Value = "  -0.38%  "
Set-TextValue "D29" "132.21"
$ws.Range("E29").Value = "  +4.08%  "
$ws.Range("D30").Value = "1.901.74"
$ws.Range("E30").Value = "  +2.08%  "
Set-TextValue "D31" "1.249"
$ws.Range("E31").Value = "  +29.11%  "
Set-TextValue "D32" "6.989"
$ws.Range("E32").Value = "  +10.86%  "
Set-TextValue "D33" "4.274"
$ws.Range("E33").Value = "  +5.91%  "
Set-TextValue "D34" "1.850"
$ws.Range("E34").Value = "  +6.46%  "
Set-TextValue "D35" "13.94"
$ws.Range("E35").Value = "  +13.58%  "
Set-TextValue "D36" "0.08846"
$ws.Range("E36").Value = "  +4.43%  "
Set-TextValue "D37" "5.637"
$ws.Range("E37").Value = "  +5.58%  "
Set-TextValue "D39" "9.188"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("E40").Value = "  +3.67%  "
Set-TextValue "D41" "0.2248"
$ws.Range("E41").Value = "  +6.70%  "
Set-TextValue "D42" "1.278"
$ws.Range("E42").Value = "  +1.25%  "
Set-TextValue "D43" "0.6484"
$ws.Range("E43").Value = "  +5.26%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D44" "0.9974"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "14.05"
$ws.Range("E45").Value = "  +8.25%  "
Set-TextValue "D46" "0.6189"
$ws.Range("E46").Value = "  +4.24%  "
Set-TextValue "D47" "3.825"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E48").Value = "  +6.00%  "
Set-TextValue "D49" "130.49"
$ws.Range("E49").Value = "  +2.70%  "
Set-TextValue "D50" "0.07332"
$ws.Range("E50").Value = "  +1.75%  "
Set-TextValue "D51" "80.09"
$ws.Range("E51").Value = "  +5.51%  "
